# Powerpoint writer: avoid extra blank lines before author.
# (In the case where there is no subtitle.)
#
# Remove the empty "Subtitle 2" placeholder shape (the one that only
# contains two line breaks) from slide 1, so the title slide has no
# subtitle shape at all.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Subtitle 2") {
        $sh.Cut()
    }
}
